# Add a new "Team Points" worksheet at the end of the workbook, matching
# the other per-category sheets (Points, Assists, Rebounds, 3PM, ...).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Team Points"

# Header row
$headers = @("Game Time (PST)", "Opponent", "Team Points", "Opponent Points", "Game Total Points")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows: Game Time (PST), Opponent, Team Points, Opponent Points, Game Total Points
$data = @(
    @("2025-10-22", "POR", 118, 114, 232),
    @("2025-10-24", "LAL", 110, 128, 238),
    @("2025-10-26", "IND", 114, 110, 224),
    @("2025-10-27", "DEN", 114, 127, 241),
    @("2025-10-29", "LAL", 115, 116, 231),
    @("2025-11-01", "CHA", 122, 105, 227),
    @("2025-11-03", "BKN", 125, 109, 234),
    @("2025-11-05", "NYK", 114, 137, 251),
    @("2025-11-07", "UTA", 137, 97, 234),
    @("2025-11-09", "SAC", 144, 117, 261),
    @("2025-11-10", "UTA", 120, 113, 233)
)

# Format the date column as text first so Excel stores the literal
# "YYYY-MM-DD" strings instead of auto-converting them to date serials.
$ws.Range("A2:A12").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    $ws.Cells.Item($excelRow, 1).Value = $row[0]
    $ws.Cells.Item($excelRow, 2).Value = $row[1]
    $ws.Cells.Item($excelRow, 3).Value = $row[2]
    $ws.Cells.Item($excelRow, 4).Value = $row[3]
    $ws.Cells.Item($excelRow, 5).Value = $row[4]
}

[void]$ws.Range("A1").Select()
